$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1): Spanish labels -> English field names ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Normalize capitalization of connector words (de/del/el/la/las/los) ---
# in state/municipality names throughout the data rows.

$ws.Range('B16').Value = 'Amatenango De La Frontera'
$ws.Range('B18').Value = 'Bejucal De Ocampo'
$ws.Range('B24').Value = 'Comitán De Domínguez'
$ws.Range('B34').Value = 'Mazapa De Madero'
$ws.Range('B39').Value = 'San Cristóbal De Las Casas'
$ws.Range('A56').Value = 'Ciudad De México'
$ws.Range('B58').Value = 'Cuajimalpa De Morelos'
$ws.Range('A70').Value = 'Coahuila De Zaragoza'
$ws.Range('A81').Value = 'Estado De México'
$ws.Range('B81').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B82').Value = 'Almoloya De Alquisiras'
$ws.Range('B83').Value = 'Almoloya De Juárez'
$ws.Range('B87').Value = 'Atizapán De Zaragoza'
$ws.Range('B99').Value = 'Ecatepec De Morelos'
$ws.Range('B102').Value = 'Ixtapan De La Sal'
$ws.Range('B108').Value = 'Naucalpan De Juárez'
$ws.Range('B112').Value = 'San Felipe Del Progreso'
$ws.Range('B113').Value = 'San Martín De Las Pirámides'
$ws.Range('B121').Value = 'Tenango Del Valle'
$ws.Range('B128').Value = 'Tlalnepantla De Baz'
$ws.Range('B131').Value = 'Valle De Bravo'
$ws.Range('B152').Value = 'Valle De Santiago'
$ws.Range('B156').Value = 'Acapulco De Juárez'
$ws.Range('B158').Value = 'Ayutla De Los Libres'
$ws.Range('B159').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B160').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B167').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B168').Value = 'Iguala De La Independencia'
$ws.Range('B177').Value = 'Taxco De Alarcón'
$ws.Range('B179').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B182').Value = 'Tlalixtaquilla De Maldonado'
$ws.Range('B183').Value = 'Técpan De Galeana'
$ws.Range('B192').Value = 'Huejutla De Reyes'
$ws.Range('B195').Value = 'Molango De Escamilla'
$ws.Range('B196').Value = 'Nopala De Villagrán'
$ws.Range('B197').Value = 'Pachuca De Soto'
$ws.Range('B199').Value = 'Santiago Tulantepec De Lugo Guerrero'
$ws.Range('B202').Value = 'Tulancingo De Bravo'
$ws.Range('B203').Value = 'Zacualtipán De Ángeles'
$ws.Range('B206').Value = 'Encarnación De Díaz'
$ws.Range('B208').Value = 'Lagos De Moreno'
$ws.Range('B211').Value = 'San Juan De Los Lagos'
$ws.Range('B212').Value = 'San Martín De Bolaños'
$ws.Range('B214').Value = 'Tepatitlán De Morelos'
$ws.Range('B215').Value = 'Unión De Tula'
$ws.Range('B217').Value = 'Zapotlán El Grande'
$ws.Range('A219').Value = 'Michoacán De Ocampo'
$ws.Range('B241').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B254').Value = 'Puente De Ixtla'
$ws.Range('B257').Value = 'Tlaltizapán De Zapata'
$ws.Range('B267').Value = 'Santa María Del Oro'
$ws.Range('B274').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B275').Value = 'Coicoyán De Las Flores'
$ws.Range('B276').Value = 'Fresnillo De Trujano'
$ws.Range('B277').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B278').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B279').Value = 'Huajuapan De León'
$ws.Range('B281').Value = 'Oaxaca De Juárez'
$ws.Range('B282').Value = 'Putla Villa De Guerrero'
$ws.Range('B300').Value = 'Santo Domingo De Morelos'
$ws.Range('B301').Value = 'Tlacolula De Matamoros'
$ws.Range('B302').Value = 'Villa De Tututepec'
$ws.Range('B303').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B325').Value = 'Huehuetlán El Chico'
$ws.Range('B331').Value = 'Palmar De Bravo'
$ws.Range('B337').Value = 'San Nicolás De Los Ranchos'
$ws.Range('B339').Value = 'San Salvador El Verde'
$ws.Range('B344').Value = 'Tepanco De López'
$ws.Range('B347').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B353').Value = 'Xochitlán De Vicente Suárez'
$ws.Range('B361').Value = 'Jalpan De Serra'
$ws.Range('B362').Value = 'Pinal De Amoles'
$ws.Range('B369').Value = 'Ciudad Del Maíz'
$ws.Range('B373').Value = 'San Ciro De Acosta'
$ws.Range('B400').Value = 'Acuamanala De Miguel Hidalgo'
$ws.Range('B404').Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range('A409').Value = 'Veracruz De Ignacio De La Llave'
$ws.Range('B411').Value = 'Amatlán De Los Reyes'
$ws.Range('B415').Value = 'Boca Del Río'
$ws.Range('B418').Value = 'Cosamaloapan De Carpio'
$ws.Range('B419').Value = 'Cosautlán De Carvajal'
$ws.Range('B427').Value = 'Ignacio De La Llave'
$ws.Range('B430').Value = 'Martínez De La Torre'
$ws.Range('B433').Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Range('B435').Value = 'Ozuluama De Mascareñas'
$ws.Range('B437').Value = 'Paso De Ovejas'
$ws.Range('B439').Value = 'Poza Rica De Hidalgo'
$ws.Range('B443').Value = 'Soledad De Doblado'
$ws.Range('A464').Value = 'Total'

# --- Remove trailing metadata/footer rows (466-470) ---
$ws.Range("A466:A470").EntireRow.Delete()
